$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sold-out "DECLOPHEN 75MG/3ML 3 AMPOULES" line item (row 8) from the
# shortage table. Deleting the entire row shifts every following row up by one,
# which Excel also takes care of for the merged-cell ranges automatically.
$ws.Rows("8:8").Delete()

# Renumber the "م" (item #) column for the rows that moved up so the sequence
# stays contiguous (1, 2, 3, 4, 5).
$ws.Range("A8").Value2 = 2
$ws.Range("A9").Value2 = 3
$ws.Range("A10").Value2 = 4
$ws.Range("A11").Value2 = 5

# Recompute the price total now that the deleted row's price (11.88) is gone:
# 148.74 - 11.88 = 136.86. The total now lives one row higher, at P12.
$ws.Range("P12").Value2 = 136.86

# Refresh the generated/printed timestamp footer (now one row higher, at A13)
# to reflect the new export time.
$ws.Range("A13").Value2 = "Tuesday, 9 September, 2025 10:15 AM"
